$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.339.84"
$ws.Range("E2").Value = "  +2.47%  "

$ws.Range("D3").Value = "2.978.33"
$ws.Range("E3").Value = "  +1.10%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'565.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.27%  "

$ws.Range("D6").Value = "'137.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.21%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("D9").Value = "2.971.58"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("E10").Value = "  +3.39%  "

$ws.Range("D11").Value = "'5.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.84%  "

$ws.Range("D12").Value = "'0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("E13").Value = "  +3.36%  "

$ws.Range("E14").Value = "  +2.26%  "

$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").Value = "3.470.32"
$ws.Range("E16").Value = "  +1.12%  "

$ws.Range("E17").Value = "  +2.15%  "

$ws.Range("D18").Value = "2.974.44"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").Value = "59.338.74"
$ws.Range("E19").Value = "  +2.47%  "

$ws.Range("D20").Value = "'435.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.43%  "

$ws.Range("D21").Value = "'13.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("D22").Value = "'0.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.17%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").Value = "'13.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.52%  "

$ws.Range("D25").Value = "'79.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").Value = "'2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.81%  "

$ws.Range("E29").Value = "  +2.24%  "

$ws.Range("D30").Value = "'7.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.19%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'25.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.87%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.105"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.65%  "

$ws.Range("D34").Value = "0.0₃0768"
$ws.Range("E34").Value = "  +9.10%  "

$ws.Range("D35").Value = "'5.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.47%  "

$ws.Range("D36").Value = "'0.984"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.97%  "

$ws.Range("D37").Value = "'2.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").Value = "'48.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").Value = "'8.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "

$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("D41").Value = "'400.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.07%  "

$ws.Range("D42").Value = "'0.0350"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("D43").Value = "2.739.12"
$ws.Range("E43").Value = "  +1.32%  "

$ws.Range("E44").Value = "  -3.25%  "

$ws.Range("E45").Value = "  +5.42%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").Value = "'34.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +19.38%  "

$ws.Range("D48").Value = "'122.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.47%  "

$ws.Range("E49").Value = "  +2.16%  "

$ws.Range("E50").Value = "  +1.23%  "

$ws.Range("D51").Value = "'23.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.20%  "

